$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.070.89'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '1.675.20'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '215.20'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +2.36%  '
$ws.Range('D9').Value = '21.25'
$ws.Range('E9').Value = '  +4.78%  '
$ws.Range('E10').Value = '  +0.10%  '
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('D12').Value = '1.912.17'
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('D13').Value = '1.676.87'
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('D15').Value = '0.535'
$ws.Range('E15').Value = '  +1.78%  '
$ws.Range('D16').Value = '66.08'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('D17').Value = '27.058.05'
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('D18').Value = '237.44'
$ws.Range('E18').Value = '  +1.85%  '
$ws.Range('E19').Value = '  +1.49%  '
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('D23').Value = '9.31'
$ws.Range('E23').Value = '  +1.77%  '
$ws.Range('E24').Value = '  -1.66%  '
$ws.Range('D25').Value = '147.48'
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('D26').Value = '7.22'
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('D27').Value = '16.33'
$ws.Range('E27').Value = '  +2.48%  '
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('D33').Value = '1.530.18'
$ws.Range('E33').Value = '  +5.03%  '
$ws.Range('E34').Value = '  +1.66%  '
$ws.Range('E35').Value = '  +3.60%  '
$ws.Range('D36').Value = '2.40'
$ws.Range('E36').Value = '  -0.75%  '
$ws.Range('D37').Value = '0.592'
$ws.Range('E37').Value = '  +1.67%  '
$ws.Range('D38').Value = '0.916'
$ws.Range('E38').Value = '  +1.85%  '
$ws.Range('E39').Value = '  +2.24%  '
$ws.Range('E40').Value = '  +2.77%  '
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('D42').Value = '67.65'
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '5.52'
$ws.Range('E43').Value = '  -3.71%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '2.26'
$ws.Range('E44').Value = '  -1.51%  '
$ws.Range('D45').Value = '1.819.64'
$ws.Range('E45').Value = '  +0.59%  '
$ws.Range('E46').Value = '  +0.28%  '
$ws.Range('D47').Value = '90.62'
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('E48').Value = '  +1.08%  '
$ws.Range('E49').Value = '  +2.21%  '
$ws.Range('D50').Value = '8.03'
$ws.Range('E50').Value = '  +4.98%  '
$ws.Range('D51').Value = '0.0510'
$ws.Range('E51').Value = '  +0.61%  '
